$wb = $excel.ActiveWorkbook

# Handback status report generation: refresh handoff/handback timestamps
# for the second source file (4b732e66-d30b-4fde-b4a3-e04f9979acbe.md)
# across the Overview sheet and each locale sheet.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-28 04:47:46"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-28 04:47:42"
$wsZhCn.Range("K3").Value = "2016-08-28 04:47:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-28 04:47:46"
$wsDeDe.Range("K3").Value = "2016-08-28 04:48:09"
